# Updated legacy GSC export data (Coverage.xlsx)

$wb = $excel.ActiveWorkbook

# Helper: write a literal text value into a cell without letting Excel's
# automatic "looks like a date" detection turn it into a date serial
# number. We do this by building the text via a formula (a quoted
# string literal always yields a text result) and then copying /
# paste-special-ing the computed value back over itself, which bakes the
# formula down into a plain, statically-typed string cell.
function Set-TextValue {
    param($cell, [string]$text)
    $cell.Formula = '="' + $text + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)  # xlPasteValues
}

# --- Sheet "Chart" -----------------------------------------------------
$chart = $wb.Worksheets.Item("Chart")

# Row 2 was a placeholder entry for 2025-11-03 with blank Not-indexed /
# Indexed counts. Remove it outright; every later day shifts up a row.
$chart.Rows.Item(2).Delete()

# Four additional days of data (2026-01-28 .. 2026-01-31) are appended
# at the bottom of the table (now rows 87-90).
Set-TextValue $chart.Cells.Item(87, 1) "2026-01-28"
$chart.Cells.Item(87, 2).Value = 106.0
$chart.Cells.Item(87, 3).Value = 29.0
$chart.Cells.Item(87, 4).Value = 0.0

Set-TextValue $chart.Cells.Item(88, 1) "2026-01-29"
$chart.Cells.Item(88, 2).Value = 106.0
$chart.Cells.Item(88, 3).Value = 29.0
$chart.Cells.Item(88, 4).Value = 4.0

Set-TextValue $chart.Cells.Item(89, 1) "2026-01-30"
$chart.Cells.Item(89, 2).Value = 106.0
$chart.Cells.Item(89, 3).Value = 29.0
$chart.Cells.Item(89, 4).Value = 6.0

Set-TextValue $chart.Cells.Item(90, 1) "2026-01-31"
$chart.Cells.Item(90, 2).Value = 106.0
$chart.Cells.Item(90, 3).Value = 29.0
$chart.Cells.Item(90, 4).Value = 0.0

# --- Sheet "Critical issues" -------------------------------------------
# Refreshed page counts for a few rows.
$critical = $wb.Worksheets.Item("Critical issues")
$critical.Cells.Item(2, 4).Value = 82.0
$critical.Cells.Item(5, 4).Value = 10.0
$critical.Cells.Item(6, 4).Value = 7.0
